$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price (D) and Volume(1h) (E) values are written with a leading
# apostrophe literal (the doubled quote inside the PS single-quoted
# string, e.g. '''0.512') so numeric-looking text (e.g. "0.512")
# is stored as text instead of being coerced into an Excel number,
# matching the original workbook (inline string / text cells).

# --- rows 2-48: updated Price / Volume(1h) figures ---
$ws.Range("D2").Value = '''41.766.40'
$ws.Range("E2").Value = '''  +0.19%  '
$ws.Range("D3").Value = '''2.478.50'
$ws.Range("E3").Value = '''  +0.54%  '
$ws.Range("E4").Value = '''  +0.08%  '
$ws.Range("D5").Value = '''321.11'
$ws.Range("E5").Value = '''  +1.97%  '
$ws.Range("D6").Value = '''92.36'
$ws.Range("E6").Value = '''  +0.09%  '
$ws.Range("E7").Value = '''  +0.47%  '
$ws.Range("E8").Value = '''  +0.05%  '
$ws.Range("D9").Value = '''0.512'
$ws.Range("E9").Value = '''  -0.03%  '
$ws.Range("E10").Value = '''  +2.35%  '
$ws.Range("D11").Value = '''33.11'
$ws.Range("E11").Value = '''  +2.04%  '
$ws.Range("E12").Value = '''  -0.66%  '
$ws.Range("D13").Value = '''2.860.59'
$ws.Range("D14").Value = '''6.90'
$ws.Range("E14").Value = '''  +0.53%  '
$ws.Range("D15").Value = '''15.56'
$ws.Range("E15").Value = '''  -1.46%  '
$ws.Range("D16").Value = '''2.479.22'
$ws.Range("E16").Value = '''  -0.06%  '
$ws.Range("D17").Value = '''0.796'
$ws.Range("E17").Value = '''  +2.27%  '
$ws.Range("D18").Value = '''41.707.11'
$ws.Range("E18").Value = '''  +0.14%  '
$ws.Range("E19").Value = '''  -0.52%  '
$ws.Range("D20").Value = '''0.0₃0944'
$ws.Range("E20").Value = '''  -0.55%  '
$ws.Range("D21").Value = '''70.72'
$ws.Range("E21").Value = '''  +0.07%  '
$ws.Range("D22").Value = '''11.27'
$ws.Range("E22").Value = '''  -1.17%  '
$ws.Range("D23").Value = '''240.02'
$ws.Range("E23").Value = '''  +0.50%  '
$ws.Range("E24").Value = '''  +1.79%  '
$ws.Range("E25").Value = '''  +1.91%  '
$ws.Range("E26").Value = '''  +0.03%  '
$ws.Range("D27").Value = '''25.02'
$ws.Range("E27").Value = '''  +2.32%  '
$ws.Range("E28").Value = '''  -0.52%  '
$ws.Range("E29").Value = '''  +0.17%  '
$ws.Range("D30").Value = '''36.71'
$ws.Range("E30").Value = '''  +4.37%  '
$ws.Range("D31").Value = '''157.48'
$ws.Range("E31").Value = '''  +1.12%  '
$ws.Range("E32").Value = '''  -0.70%  '
$ws.Range("D34").Value = '''0.0766'
$ws.Range("E34").Value = '''  +0.67%  '
$ws.Range("E35").Value = '''  -0.60%  '
$ws.Range("D36").Value = '''17.22'
$ws.Range("E36").Value = '''  -1.27%  '
$ws.Range("E37").Value = '''  +3.26%  '
$ws.Range("E38").Value = '''  +1.46%  '
$ws.Range("D39").Value = '''2.88'
$ws.Range("E39").Value = '''  -0.29%  '
$ws.Range("E40").Value = '''  +0.95%  '
$ws.Range("D41").Value = '''4.01'
$ws.Range("E41").Value = '''  +1.80%  '
$ws.Range("E42").Value = '''  -1.99%  '
$ws.Range("D43").Value = '''1.997.25'
$ws.Range("E43").Value = '''  +1.14%  '
$ws.Range("E44").Value = '''  +0.73%  '
$ws.Range("E45").Value = '''  -0.23%  '
$ws.Range("D46").Value = '''2.98'
$ws.Range("E46").Value = '''  +1.35%  '
$ws.Range("D47").Value = '''9.45'
$ws.Range("E47").Value = '''  +5.01%  '
$ws.Range("D48").Value = '''2.738.13'
$ws.Range("E48").Value = '''  +1.42%  '

# --- rows 49-50: Aave and BitcoinSV swapped ranking positions ---
$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D49").Value = '''76.22'
$ws.Range("E49").Value = '''  +5.49%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '''97.76'
$ws.Range("E50").Value = '''  +1.01%  '
$ws.Range("D51").Value = '''67.43'
$ws.Range("E51").Value = '''  +0.73%  '
